$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "house_code" (鸽棚号) column D is being dropped from the template.
# Deleting the entire column shifts every cell/column to its right one
# position to the left (D disappears, E->D, F->E, ... L->K), which is
# exactly the row/col/merge/dimension shift seen in the diff.
$ws.Columns("D:D").Delete()

# Update the instructional note (now living in F2 after the column
# shift) to drop the house_code explanation line and reword
# "其余3项" -> "其余项目".
$note = $ws.Range("F2")
$note.Value2 = "环号、团长和玩家为必填项，其余项目可选项。`n标签：如果有多个标签请用逗号（,）分割，`n文件中的标签如果不存在，会自动创建`n"

$headline = $note.Characters(1, 21)
$headline.Font.Name = "等线 (正文)"
$headline.Font.Size = 14
$headline.Font.Bold = $true
$headline.Font.Color = 255

$body = $note.Characters(22, 40)
$body.Font.Name = "等线"
$body.Font.Size = 11

# Re-typing the note re-triggers autofit on the row; put the height back
# to the sheet's normal automatic height (no stored override).
$note.EntireRow.AutoFit()

# Restore the selection to what it would be after this edit.
$null = $ws.Range("F2:K20").Select()
